$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the HED tag text in D2 (remove "Organizational/ExampleFileOrg, ")
$ws.Range("D2").Value = "Event/Sensory event, (Participant/Cognitive state/Awake ~ Participant/Trait/Age/15 ~ Item/Sound/Siren, Item/Object/Manmade/Vehicle, Attribute/Sensory/Visual/Color/RGB color/RGB Red/100), (Definition/ExampleFileDef, (Item/ExampleFileTag))"

# Update the HED tag text in D4 (rename "Label-def/ExampleFileDef" to "def/ExampleFileDef")
$ws.Range("D4").Value = "Event/Sensory event, (Participant/Cognitive state/Awake ~ Participant/Trait/Age/15 ~ Item/Sound/Siren, Item/Object/Manmade/Vehicle, Attribute/Sensory/Visual/Color/RGB color/RGB Red/100), def/ExampleFileDef"

# Update the HED tag text in D5 (rename "Label-def/ExampleDef" to "def/ExampleDef")
$ws.Range("D5").Value = "Event/Sensory event, (Participant/Cognitive state/Awake ~ Participant/Trait/Age/15 ~ Item/Sound/Siren, Item/Object/Manmade/Vehicle, Attribute/Sensory/Visual/Color/RGB color/RGB Red/100), def/ExampleDef"

# Move the active cell selection from D5 to D4
$ws.Range("D4").Select()

# Adjust row 2 height
$ws.Rows.Item(2).RowHeight = 82.05
